$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rename to snake_case identifiers ---
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# --- Capitalize connector words (de/del/la/las/los/el/y) in place/municipality names ---
$renames = @(
    @('B6', 'Rincón De Romos'),
    @('B26', 'Amatenango De La Frontera'),
    @('B29', 'Bejucal De Ocampo'),
    @('B31', 'Benemérito De Las Américas'),
    @('B37', 'Chiapa De Corzo'),
    @('B41', 'Comitán De Domínguez'),
    @('B57', 'Mazapa De Madero'),
    @('B63', 'Ocozocoautla De Espinosa'),
    @('B68', 'Salto De Agua'),
    @('B69', 'San Cristóbal De Las Casas'),
    @('B92', 'Coyame Del Sotol'),
    @('B94', 'Guadalupe Y Calvo'),
    @('B95', 'Hidalgo Del Parral'),
    @('B113', 'San Juan De Sabinas'),
    @('A123', 'Ciudad De México'),
    @('B127', 'Cuajimalpa De Morelos'),
    @('B148', 'Nombre De Dios'),
    @('A156', 'Estado De México'),
    @('B156', 'Acambay De Ruíz Castañeda'),
    @('B158', 'Almoloya De Alquisiras'),
    @('B159', 'Almoloya De Juárez'),
    @('B162', 'Atizapán De Zaragoza'),
    @('B166', 'Chapa De Mota'),
    @('B170', 'Coacalco De Berriozábal'),
    @('B174', 'Ecatepec De Morelos'),
    @('B177', 'Ixtapan De La Sal'),
    @('B178', 'Ixtapan Del Oro'),
    @('B190', 'Naucalpan De Juárez'),
    @('B198', 'San Antonio La Isla'),
    @('B199', 'San Felipe Del Progreso'),
    @('B207', 'Tenango Del Valle'),
    @('B213', 'Tlalnepantla De Baz'),
    @('B217', 'Valle De Bravo'),
    @('B218', 'Valle De Chalco Solidaridad'),
    @('B219', 'Villa De Allende'),
    @('B228', 'San Miguel De Allende'),
    @('B229', 'Apaseo El Alto'),
    @('B230', 'Apaseo El Grande'),
    @('B237', 'Dolores Hidalgo Cuna De La Independencia Nacional'),
    @('B240', 'Jaral Del Progreso'),
    @('B250', 'San Diego De La Unión'),
    @('B252', 'San Francisco Del Rincón'),
    @('B254', 'San Luis De La Paz'),
    @('B256', 'Santa Cruz De Juventino Rosas'),
    @('B257', 'Silao De La Victoria'),
    @('B261', 'Valle De Santiago'),
    @('B267', 'Acapulco De Juárez'),
    @('B270', 'Ajuchitlán Del Progreso'),
    @('B271', 'Alcozauca De Guerrero'),
    @('B274', 'Atenango Del Río'),
    @('B276', 'Atoyac De Álvarez'),
    @('B277', 'Ayutla De Los Libres'),
    @('B280', 'Buenavista De Cuéllar'),
    @('B281', 'Chilapa De Álvarez'),
    @('B282', 'Chilpancingo De Los Bravo'),
    @('B286', 'Coyuca De Benítez'),
    @('B287', 'Coyuca De Catalán'),
    @('B290', 'Cuetzala Del Progreso'),
    @('B291', 'Cutzamala De Pinzón'),
    @('B296', 'Huitzuco De Los Figueroa'),
    @('B297', 'Iguala De La Independencia'),
    @('B299', 'Ixcateopan De Cuauhtémoc'),
    @('B300', 'Zihuatanejo De Azueta'),
    @('B302', 'La Unión De Isidoro Montes De Oca'),
    @('B305', 'Mártir De Cuilapan'),
    @('B314', 'Taxco De Alarcón'),
    @('B316', 'Técpan De Galeana'),
    @('B318', 'Tepecoacuilco De Trujano'),
    @('B320', 'Tixtla De Guerrero'),
    @('B323', 'Tlalixtaquilla De Maldonado'),
    @('B324', 'Tlapa De Comonfort'),
    @('B339', 'Atotonilco El Grande'),
    @('B344', 'Cuautepec De Hinojosa'),
    @('B346', 'Huasca De Ocampo'),
    @('B348', 'Huejutla De Reyes'),
    @('B351', 'Jacala De Ledezma'),
    @('B357', 'Mineral Del Monte'),
    @('B358', 'Mixquiahuala De Juárez'),
    @('B359', 'Molango De Escamilla'),
    @('B361', 'Nopala De Villagrán'),
    @('B362', 'Omitlán De Juárez'),
    @('B363', 'Pachuca De Soto'),
    @('B366', 'Progreso De Obregón'),
    @('B371', 'Santiago De Anaya'),
    @('B375', 'Tenango De Doria'),
    @('B376', 'Tepehuacán De Guerrero'),
    @('B377', 'Tepeji Del Río De Ocampo'),
    @('B378', 'Tezontepec De Aldama'),
    @('B384', 'Tula De Allende'),
    @('B385', 'Tulancingo De Bravo'),
    @('B388', 'Zacualtipán De Ángeles'),
    @('B392', 'Acatlán De Juárez'),
    @('B395', 'Autlán De Navarro'),
    @('B401', 'Concepción De Buenos Aires'),
    @('B405', 'Encarnación De Díaz'),
    @('B409', 'Ixtlahuacán Del Río'),
    @('B412', 'Jilotlán De Los Dolores'),
    @('B416', 'Lagos De Moreno'),
    @('B423', 'San Juan De Los Lagos'),
    @('B426', 'San Miguel El Alto'),
    @('B428', 'Talpa De Allende'),
    @('B429', 'Tamazula De Gordiano'),
    @('B433', 'Tepatitlán De Morelos'),
    @('B435', 'Tizapán El Alto'),
    @('B436', 'Tlajomulco De Zúñiga'),
    @('B438', 'Unión De Tula'),
    @('B440', 'Zacoalco De Torres'),
    @('B442', 'Zapotlán El Grande'),
    @('B495', 'Tiquicheo De Nicolás Romero'),
    @('B512', 'Coatlán Del Río'),
    @('B518', 'Jonacatepec De Leandro Valle'),
    @('B522', 'Puente De Ixtla'),
    @('B526', 'Tetela Del Volcán'),
    @('B528', 'Tlaltizapán De Zapata'),
    @('B538', 'Bahía De Banderas'),
    @('B541', 'Ixtlán Del Río'),
    @('B545', 'Santa María Del Oro'),
    @('B557', 'San Nicolás De Los Garza'),
    @('B561', 'Acatlán De Pérez Figueroa'),
    @('B567', 'Coicoyán De Las Flores'),
    @('B571', 'Cuilápam De Guerrero'),
    @('B572', 'Heroica Ciudad De Ejutla De Crespo'),
    @('B573', 'Heroica Ciudad De Tlaxiaco'),
    @('B574', 'Huautla De Jiménez'),
    @('B575', 'Ixtlán De Juárez'),
    @('B576', 'Heroica Ciudad De Juchitán De Zaragoza'),
    @('B583', 'Mártires De Tacubaya'),
    @('B585', 'Miahuatlán De Porfirio Díaz'),
    @('B588', 'Oaxaca De Juárez'),
    @('B589', 'Ocotlán De Morelos'),
    @('B590', 'Pinotepa De Don Luis'),
    @('B592', 'Putla Villa De Guerrero'),
    @('B604', 'San Dionisio Del Mar'),
    @('B610', 'San Francisco Del Mar'),
    @('B618', 'San José Del Progreso'),
    @('B621', 'San Juan Bautista Lo De Soto'),
    @('B630', 'San Juan De Los Cués'),
    @('B650', 'San Miguel El Grande'),
    @('B658', 'San Pablo Villa De Mitla'),
    @('B660', 'San Pedro El Alto'),
    @('B672', 'San Pedro Y San Pablo Teposcolula'),
    @('B682', 'Santa Ana Del Valle'),
    @('B721', 'Santiago Del Río'),
    @('B739', 'Santo Domingo De Morelos'),
    @('B748', 'Tamazulápam Del Espíritu Santo'),
    @('B749', 'Tataltepec De Valdés'),
    @('B750', 'Tezoatlán De Segura Y Luna'),
    @('B751', 'Tlacolula De Matamoros'),
    @('B752', 'Tlalixtac De Cabrera'),
    @('B753', 'Villa De Tututepec De Melchor Ocampo'),
    @('B754', 'Villa De Zaachila'),
    @('B755', 'Villa Sola De Vega'),
    @('B756', 'Zimatlán De Álvarez'),
    @('B777', 'Cuayuca De Andrade'),
    @('B787', 'Huehuetlán El Chico'),
    @('B791', 'Izúcar De Matamoros'),
    @('B795', 'Los Reyes De Juárez'),
    @('B799', 'Palmar De Bravo'),
    @('B805', 'San Salvador El Seco'),
    @('B806', 'San Salvador El Verde'),
    @('B812', 'Tepatlaxco De Hidalgo'),
    @('B814', 'Tepexi De Rodríguez'),
    @('B818', 'Tlacotepec De Benito Juárez'),
    @('B828', 'Xayacatlán De Bravo'),
    @('B844', 'Amealco De Bonfil'),
    @('B846', 'Cadereyta De Montes'),
    @('B849', 'Jalpan De Serra'),
    @('B850', 'Landa De Matamoros'),
    @('B852', 'Pinal De Amoles'),
    @('B854', 'San Juan Del Río'),
    @('B865', 'Armadillo De Los Infante'),
    @('B866', 'Axtla De Terrazas'),
    @('B870', 'Ciudad Del Maíz'),
    @('B878', 'Mexquitic De Carmona'),
    @('B884', 'San Ciro De Acosta'),
    @('B888', 'Santa María Del Río'),
    @('B889', 'Soledad De Graciano Sánchez'),
    @('B897', 'Villa De La Paz'),
    @('B898', 'Villa De Ramos'),
    @('B899', 'Villa De Reyes'),
    @('B925', 'Jalpa De Méndez'),
    @('B950', 'Soto La Marina'),
    @('B959', 'Nanacamilpa De Mariano Arista'),
    @('B960', 'San Pablo Del Monte'),
    @('B963', 'Tetla De La Solidaridad'),
    @('B973', 'Alto Lucero De Gutiérrez Barrios'),
    @('B977', 'Amatlán De Los Reyes'),
    @('B986', 'Boca Del Río'),
    @('B988', 'Camarón De Tejeda'),
    @('B992', 'Castillo De Teayo'),
    @('B994', 'Cazones De Herrera'),
    @('B1007', 'Cosamaloapan De Carpio'),
    @('B1008', 'Cosautlán De Carvajal'),
    @('B1024', 'Hueyapan De Ocampo'),
    @('B1025', 'Ignacio De La Llave'),
    @('B1028', 'Ixhuacán De Los Reyes'),
    @('B1029', 'Ixhuatlán De Madero'),
    @('B1030', 'Ixhuatlán Del Café'),
    @('B1039', 'Las Vigas De Ramírez'),
    @('B1040', 'Lerdo De Tejada'),
    @('B1044', 'Martínez De La Torre'),
    @('B1054', 'Ozuluama De Mascareñas'),
    @('B1057', 'Paso De Ovejas'),
    @('B1058', 'Paso Del Macho'),
    @('B1061', 'Poza Rica De Hidalgo'),
    @('B1069', 'Sayula De Alemán'),
    @('B1072', 'Soledad De Doblado'),
    @('B1076', 'Tatahuicapan De Juárez'),
    @('B1095', 'Tlacotepec De Mejía'),
    @('B1106', 'Vega De Alatorre'),
    @('B1120', 'Concepción Del Oro')
)

foreach ($pair in $renames) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# --- Fix tiny floating point rounding on two percentage cells ---
$ws.Range('D303').Value = 0.009046418508082457
$ws.Range('D390').Value = 0.09268871422215633

# --- Remove trailing footnote/metadata rows (sample size, source, author, date) ---
$ws.Range('A1143:D1147').EntireRow.Delete()

Write-Host "edit complete"
